$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = @([double]"1.665374230732391e-22", [double]"2.524364338824382e-32", [double]"1.352704871469615e-22", [double]"-1.352749823090494e-22", [double]"-2.760134611073494e-33", [double]"-5.128843399645042e-33", [double]"5.128824909122333e-33", [double]"4.344454514954092e-22", [double]"7.093829170853745e-22", [double]"-4.992583390432867e-22", [double]"7.199091411446751e-22", [double]"-1.20714640863934e-11", [double]"-1.287622835881957e-11", [double]"-1.379595895587827e-11", [double]"-1.48571865678687e-11", 1)
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row[$i] }

$row = @([double]"1.379595441334928e-11", [double]"2.101680114068488e-21", [double]"1.084597122394622e-17", [double]"-1.084389355574739e-17", [double]"-1.423503335560283e-20", [double]"1.143239369995058e-21", [double]"-1.177489229829343e-21", [double]"2.759190919626479e-11", [double]"1.273014087395436e-11", [double]"-3.828693187204522e-11", [double]"1.75902750105536e-17", -1.0, [double]"3.366380189682996e-15", [double]"-2.414850701961177e-15", [double]"3.249067041065305e-16", [double]"-1.207132318112019e-11")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row[$i] }

$row = @([double]"6.569568230479623e-18", [double]"6.428040652633421e-23", [double]"9.10471538973065e-12", [double]"-9.105017948701963e-12", [double]"7.105090464649326e-21", [double]"-1.627035819208651e-21", [double]"1.555382260861803e-21", [double]"-2.295471486847637e-17", [double]"7.956372559978297e-12", [double]"-2.392933034335652e-11", [double]"3.219057925514334e-11", [double]"-2.604846392911482e-16", [double]"-9.782496747449399e-16", [double]"-3.272359287040641e-15", -1.0, [double]"-1.485749474322401e-11")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row[$i] }

$row = @([double]"-1.375888438987634e-16", [double]"1.960508313104706e-17", [double]"1.951082183501514e-11", [double]"1.951007645486497e-11", [double]"3.219040585236756e-11", [double]"8.305575431744415e-17", [double]"-8.307143615362624e-17", [double]"-2.650361544865463e-17", 0.3295527828388244, -0.9911558133533193, [double]"5.409388983293431e-16", [double]"3.862860021013119e-11", [double]"7.016479167557147e-17", [double]"8.499726758271735e-17", [double]"2.414285421803879e-11", [double]"2.476730230227017e-22")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row[$i] }

$row = @([double]"1.20714464047233e-11", [double]"1.75793882695788e-21", [double]"1.51460922534458e-17", [double]"-1.514509750148364e-17", [double]"-2.692349409069774e-20", [double]"3.701622416185696e-21", [double]"-3.592567666980135e-21", [double]"3.719393792347411e-17", [double]"2.605052244088411e-11", [double]"3.661431873046316e-12", [double]"3.862868595054531e-11", [double]"1.591012737683273e-15", [double]"2.207624148062156e-15", -1.0, [double]"4.555383847910845e-15", [double]"-1.379419262672943e-11")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $row[$i] }

$row = @(1.0, [double]"1.931434288758194e-10", [double]"9.97025833586315e-16", [double]"-9.97922205680996e-16", [double]"-6.438120148163163e-11", [double]"-2.757875719321384e-18", [double]"2.775472420211469e-18", [double]"4.701484335048023e-17", [double]"-1.150478492675772e-16", [double]"7.026111492536771e-17", [double]"-2.914335447595672e-16", [double]"1.379595895587471e-11", [double]"1.037542581496375e-23", [double]"1.207146408638367e-11", [double]"2.909823889434485e-24", [double]"1.665141958277148e-22")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $row[$i] }

$row = @([double]"-4.535232651725648e-16", [double]"9.541095414789774e-18", [double]"2.841197976332841e-16", [double]"-2.912646725097674e-16", [double]"2.414295406849748e-11", [double]"1.951046635613485e-11", [double]"-1.951039605059604e-11", [double]"1.095736445465087e-16", [double]"-3.888909814843483e-16", [double]"2.012299220649949e-16", 1, [double]"-1.352587394875078e-24", [double]"2.952979972853078e-24", [double]"3.862868507645693e-11", [double]"3.219057089704932e-11", [double]"2.912136469209888e-22")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $row[$i] }

$row = @([double]"6.438117675880091e-11", [double]"1.731020118964983e-17", [double]"5.557573719739602e-16", [double]"9.367003398047486e-16", 1.0, [double]"-1.365724286741246e-10", [double]"-1.365736229829842e-10", [double]"-5.619536286739466e-17", [double]"-1.060869609552739e-11", [double]"3.190581061690605e-11", [double]"-2.414300810704911e-11", [double]"-9.206094933747034e-22", [double]"2.234580334567129e-25", [double]"-7.534202529190264e-22", [double]"-7.771353461079868e-22", [double]"-5.04540302465266e-33")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $row[$i] }

$row = @([double]"1.565463539576653e-17", [double]"-1.336428445804674e-22", [double]"1.050541947480566e-11", [double]"-1.050576858222879e-11", [double]"3.366758073665463e-20", [double]"-3.908243614130491e-21", [double]"3.905089793985532e-21", [double]"2.414292264905585e-11", [double]"3.039238869299804e-11", [double]"4.271723893575094e-12", [double]"-2.682245524260773e-17", [double]"-2.672328160745508e-15", -1.0, [double]"-1.945501226067442e-15", [double]"1.247757597534541e-15", [double]"-1.287568532257792e-11")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $row[$i] }

$row = @([double]"-1.90610959539296e-16", [double]"-3.219056760967205e-11", [double]"2.731489413278035e-11", [double]"2.731431281809715e-11", [double]"1.326357570909331e-16", [double]"-8.20479230976462e-18", [double]"8.190220512194377e-18", 1, [double]"3.755887814602189e-16", [double]"-4.105722334831291e-16", [double]"2.69540548887559e-16", [double]"2.759191791175844e-11", [double]"2.414292817278058e-11", [double]"6.257280680690587e-24", [double]"-2.104002919209558e-24", [double]"2.094822528617904e-22")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $row[$i] }

$row = @([double]"-1.262149164225921e-15", [double]"1.154200975251423e-15", 0.7070950659087049, -0.7071185633971495, [double]"3.204114339881421e-16", [double]"-6.828662951612312e-11", [double]"6.828638341608807e-11", [double]"2.344009897043123e-16", [double]"9.469222565548587e-17", [double]"1.190172222066886e-17", [double]"2.042358407885334e-16", [double]"1.231769214506349e-24", [double]"1.485718656786489e-11", [double]"4.340275290830633e-24", [double]"1.287622835881807e-11", [double]"1.913004272763913e-22")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $row[$i] }

$row = @([double]"-2.326062693822745e-25", [double]"-1.241579225390183e-15", 0.7071184962702959, 0.7070949987796202, [double]"-6.808451113457802e-16", [double]"-6.828660696150925e-11", [double]"-6.828641886836668e-11", [double]"-3.862868507645601e-11", [double]"-9.09296605056293e-12", [double]"2.734787184700931e-11", [double]"-8.715238596914975e-26", [double]"-1.065837667222179e-21", [double]"-5.794339357851038e-22", [double]"-6.981834875761937e-28", [double]"-4.766661924972722e-22", [double]"-4.472847828317725e-33")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $row[$i] }

$row = @([double]"-1.5756342512065e-16", [double]"-2.414290007734553e-11", [double]"-8.041775851746166e-17", [double]"7.878977529051663e-17", [double]"-1.038105732383381e-17", [double]"2.731455399390334e-11", [double]"-2.731445555949354e-11", [double]"-7.29157315092532e-18", 0.944137152813714, 0.1327032541270937, [double]"2.209214383758155e-16", [double]"-1.09590350685046e-25", [double]"3.219057089704948e-11", [double]"2.759191791175612e-11", [double]"-1.675666780433509e-25", [double]"2.452460505825282e-22")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $row[$i] }

$row = @([double]"1.931434253822951e-10", -1.0, [double]"-5.797682360548315e-16", [double]"-2.271909548600767e-15", [double]"1.457283827394942e-17", [double]"-4.55242756308276e-11", [double]"-4.552440825531886e-11", [double]"-3.219057089706281e-11", [double]"-2.279430734335372e-11", [double]"-3.203745108759551e-12", [double]"3.328109496085297e-25", [double]"-6.149074082706975e-22", [double]"-7.77174641007935e-22", [double]"-5.107147640916535e-22", [double]"9.948612767191988e-35", [double]"-4.13683959236623e-33")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $row[$i] }

$row = @([double]"-9.110429739083297e-26", [double]"1.589774002716592e-16", [double]"6.828548122608369e-11", [double]"-6.828755107045754e-11", [double]"1.126734025035226e-16", 0.70710805536802, -0.7071055070027789, [double]"3.743438718584211e-26", [double]"-3.647085296111135e-11", [double]"-5.126282407211175e-12", [double]"-2.759191791175633e-11", [double]"3.431099875263232e-34", [double]"-8.260877303035707e-22", [double]"-1.065839507667986e-21", [double]"-6.695658445613472e-22", [double]"-5.4050651607276e-33")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $row[$i] }

$row = @([double]"1.243479557629784e-20", [double]"-6.438114179409828e-11", [double]"6.828764426041128e-11", [double]"6.828537506749272e-11", [double]"1.93143425382295e-10", 0.707105507002779, 0.7071080553680201, [double]"4.559432667379332e-21", [double]"2.897928656279906e-21", [double]"-2.227570613316149e-21", [double]"1.42746543255417e-21", [double]"5.024271810181825e-32", [double]"3.832215787878858e-32", [double]"3.750088912921917e-32", [double]"2.857597990731825e-32", [double]"1.415430784746575e-43")
for ($i = 0; $i -lt $row.Length; $i++) { $ws.Cells.Item(17, $i + 1).Value = $row[$i] }

Write-Host "Applied eigenvector updates"